# Applies the "row rotation" update described by the diff:
# rows 4..10 of the sheet effectively shift up by one (old row 5 -> row 4,
# old row 6 -> row 5, ..., old row 10 -> row 9) while old row 4's data
# wraps around into row 10. Columns A, B, D, E, F, G, H, Q, R are the ones
# that carry the per-row data; the rest (C, I, P, S, T, U, V, W, Y, Z, AA,
# AB, AD, AE, AG, AT, AW, AX, AY) stay constant for all these rows.
#
# Additionally, the source row that ends up on row 9 (old row 10, Id
# 102089632) had a handful of blank placeholder cells in columns J, K, N
# and AF that travel with it, so those placeholders move from row 10 to
# row 9 as part of the rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A, B, D, E, F, G, H, Q, R for rows 4 through 10.
$newValues = @{
    4  = @{ A = 102089519; B = 78098; D = "NT"; E = 6453; F = "Vedskivlav";          G = "Hertelidea botryosa";        H = "(Fr.) Printzen & Kantvilas";        Q = 436719.9989723715; R = 6826833.746279179 }
    5  = @{ A = 102089499; B = 78098; D = "NT"; E = 6453; F = "Vedskivlav";          G = "Hertelidea botryosa";        H = "(Fr.) Printzen & Kantvilas";        Q = 436795.9152835784; R = 6826835.694220046 }
    6  = @{ A = 102089518; B = 77258; D = "NT"; E = 6446; F = "Kolflarnlav";         G = "Carbonicola anthracophila";  H = "(Nyl.) Bendiksby & Timdal";         Q = 436745.1602881325; R = 6826800.40687584 }
    7  = @{ A = 102089652; B = 89633; D = "VU"; E = 65;   F = "Fläckporing";         G = "Anthoporia albobrunnea";     H = "(Romell) Karasiński & Niemelä";     Q = 436792.8352663805; R = 6826823.837742299 }
    8  = @{ A = 102089546; B = 77258; D = "NT"; E = 6446; F = "Kolflarnlav";         G = "Carbonicola anthracophila";  H = "(Nyl.) Bendiksby & Timdal";         Q = 436728.3378123537; R = 6826872.190838255 }
    9  = @{ A = 102089632; B = 76909; D = "NT"; E = 6437; F = "Blanksvart spiklav";  G = "Calicium denigratum";        H = "(Vain.) Tibell";                    Q = 436882.212493244;  R = 6826778.843434816 }
    10 = @{ A = 102089658; B = 78098; D = "NT"; E = 6453; F = "Vedskivlav";          G = "Hertelidea botryosa";        H = "(Fr.) Printzen & Kantvilas";        Q = 436500.5104121323; R = 6827329.967810398 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}

# Move the blank placeholder cells from row 10 to row 9 (mirrors the row
# that now lives on row 9 carrying its original blank J/K/N/AF cells).
$ws.Range("I10").Copy($ws.Range("J9"))
$ws.Range("I10").Copy($ws.Range("K9"))
$ws.Range("I10").Copy($ws.Range("N9"))
$ws.Range("I10").Copy($ws.Range("AF9"))

$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("AF10").ClearContents()
